# Program now automatically tracks all scores from game and writes to
# excel upon exiting. Each worksheet (one per game "scene") previously
# held every recorded visit; the tracker now only keeps what is actually
# present for that scene, so most sheets lose rows and some sheets end
# up with no visits recorded at all.

$wb = $excel.ActiveWorkbook

# --- s0 : 3 visits -> 2 visits, new timings -----------------------------
$ws = $wb.Worksheets.Item("s0")
$ws.Rows.Item(4).Delete()
$ws.Range("B2").Value = 1.059
$ws.Range("B3").Value = 1.223
$ws.Range("D3").Value = 1.141

# --- s1 : 3 visits -> 2 visits, new timings -----------------------------
$ws = $wb.Worksheets.Item("s1")
$ws.Rows.Item(4).Delete()
$ws.Range("B2").Value = 0.861
$ws.Range("B3").Value = 1.522
$ws.Range("D3").Value = 1.1915

# --- s2 : 4 visits -> 1 visit, new timing -------------------------------
$ws = $wb.Worksheets.Item("s2")
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()
$ws.Range("A3:B3").ClearContents()
$ws.Range("B2").Value = 2.276
$ws.Range("D3").Value = 2.276

# --- s5 : 1 visit -> 0 visits (no data recorded this run) ---------------
$ws = $wb.Worksheets.Item("s5")
$ws.Rows.Item(3).Delete()
$ws.Rows.Item(2).Delete()

# --- s6 : 1 visit -> 1 visit, new timing --------------------------------
$ws = $wb.Worksheets.Item("s6")
$ws.Range("B2").Value = 6.058
$ws.Range("D3").Value = 6.058

# --- s7 : 1 visit -> 0 visits (no data recorded this run) ---------------
$ws = $wb.Worksheets.Item("s7")
$ws.Rows.Item(3).Delete()
$ws.Rows.Item(2).Delete()

# --- s8 : 1 visit -> 0 visits (no data recorded this run) ---------------
$ws = $wb.Worksheets.Item("s8")
$ws.Rows.Item(3).Delete()
$ws.Rows.Item(2).Delete()
